$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Hüseyin"
$ws.Range("B1").Value = "Ünalan"
$ws.Range("C1").Value = 1515

$ws.Range("C1").Select()
